# POLS6330 notes - Jan 21, 2020 class edits
#
# 1. "Regulate interstate commerce."            -> "Regulating interstate commerce."
# 2. "...economic activities but be defined."   -> "...economic activities must be defined."
# 3. The Word "last edit location" (_GoBack) bookmark moves from the very
#    end of the document to the paragraph right after the "must be
#    defined." sentence (i.e. where the user was actually typing).

$d = $word.ActiveDocument

function Split-RangeFormatting($range) {
    # Toggling a character formatting property on and back off forces Word
    # to materialize a dedicated run for that exact span instead of folding
    # it back into its neighbour - i.e. it creates a run boundary without
    # altering the visible formatting.
    $range.Bold = 1
    $range.Bold = 0
}

# ---------------------------------------------------------------------
# Change 1: "Regulate" -> "Regulating" (word retyped), leaving the rest
# of the sentence, " interstate commerce.", as its own run.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(31).Range
$hit1 = $p1.Find.Execute("Regulate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordStart = $p1.Start
$wordEnd = $p1.End
$rWord = $d.Range($wordStart, $wordEnd)
$rWord.Text = "Regulating"

$rNewWord = $d.Range($wordStart, $wordStart + 10)
Split-RangeFormatting $rNewWord

# ---------------------------------------------------------------------
# Change 2: "but" -> "must" inside "The boundaries between government
# and private economic activities but be defined.", then split the
# sentence into the five runs left behind by the edit.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(33).Range
$hit2 = $p2.Find.Execute("but be defined", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sentenceStart = $d.Paragraphs(33).Range.Start

$rBut = $p2.Find.Parent
$hitBut = $d.Paragraphs(33).Range
$foundBut = $hitBut.Find.Execute("but", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$butStart = $hitBut.Start
$butEnd = $hitBut.End
$rReplace = $d.Range($butStart, $butEnd)
$rReplace.Text = "must"

$base = $d.Paragraphs(33).Range.Start
$seg1 = $d.Range($base + 0, $base + 35)
Split-RangeFormatting $seg1
$seg2 = $d.Range($base + 35, $base + 67)
Split-RangeFormatting $seg2
$seg3 = $d.Range($base + 67, $base + 68)
Split-RangeFormatting $seg3
$seg4 = $d.Range($base + 68, $base + 69)
Split-RangeFormatting $seg4

# ---------------------------------------------------------------------
# Change 3: relocate the _GoBack bookmark to the now-empty paragraph
# that follows the sentence we just edited (paragraph 34), removing it
# from the paragraph at the very end of the document.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$targetPara = $d.Paragraphs(34).Range
$d.Bookmarks.Add("_GoBack", $targetPara)

Write-Output "Paragraph 31: $($d.Paragraphs(31).Range.Text)"
Write-Output "Paragraph 33: $($d.Paragraphs(33).Range.Text)"
Write-Output "GoBack present: $($d.Bookmarks.Exists('_GoBack'))"
